$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value (as text, matching original inlineStr text cells)
$updates = @{
    'D2' = '245.21'
    'E2' = '-0.55%'
    'D3' = '26.94'
    'E3' = '2.28%'
    'D4' = '5.082'
    'E4' = '-0.07%'
    'D5' = '0.05698'
    'E5' = '1.69%'
    'D6' = '6.489'
    'E6' = '0.10%'
    'D8' = '0.9030'
    'E8' = '6.60%'
    'D9' = '0.1327'
    'E9' = '-0.96%'
    'D10' = '0.06879'
    'E10' = '-1.83%'
    'E11' = '-0.70%'
    'E12' = '0.07%'
    'D13' = '0.001522'
    'E13' = '0.57%'
    'D14' = '0.04106'
    'E14' = '-11.69%'
    'B15' = 'One'
    'C15' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D15' = '0.0005996'
    'E15' = '-93.95%'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.006097'
    'E16' = '-1.61%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.507'
    'E17' = '-2.81%'
    'B18' = 'GateToken'
    'C18' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D18' = '3.007'
    'E18' = '-0.18%'
    'B19' = 'BTSEToken'
    'C19' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D19' = '2.316'
    'E19' = '12.69%'
    'D21' = '0.03196'
    'E21' = '0.26%'
    'E22' = '-1.79%'
    'D23' = '3.553'
    'E23' = '-5.33%'
    'E24' = '1.70%'
    'D25' = '0.001215'
    'E25' = '-2.67%'
    'D26' = '0.003971'
    'E26' = '-13.46%'
    'D27' = '0.00009895'
    'E27' = '3.04%'
    'E28' = '-25.30%'
    'D40' = '0.03696'
    'E40' = '0.67%'
    'D41' = '0.005717'
    'E41' = '-7.37%'
    'E42' = '0.05%'
    'D43' = '0.002367'
    'E43' = '-5.33%'
    'D44' = '0.009389'
    'E44' = '7.20%'
    'D45' = '0.00005202'
    'E45' = '-1.80%'
    'E46' = '-0.08%'
    'E47' = '-0.09%'
    'D48' = '0.002575'
    'E48' = '-4.36%'
    'E49' = '-0.08%'
    'E50' = '-0.08%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}

Write-Host "Applied $($updates.Count) cell updates"